# "Raw and Clean Data from SSA for August 12th"
# Appends the August 12, 2020 (serial 44055) observation row to each of the
# daily tracking sheets, fills in the corresponding BV (day 74) column on
# control_obs, extends its running-total formula, and leaves the workbook
# focused on control_obs the way the author left it.

$wb = $excel.ActiveWorkbook

$wsOut     = $wb.Worksheets.Item("out_vars")
$wsDx      = $wb.Worksheets.Item("dates_dx")
$wsSx      = $wb.Worksheets.Item("dates_sx")
$wsDeaths  = $wb.Worksheets.Item("dates_deaths")
$wsControl = $wb.Worksheets.Item("control_obs")

# ---------------------------------------------------------------------------
# out_vars: new row 74 (copy formats forward from row 73, then overwrite values)
# ---------------------------------------------------------------------------
$wsOut.Range("A73:J73").Copy() | Out-Null
$wsOut.Range("A74:J74").PasteSpecial(-4122) | Out-Null
$wsOut.Application.CutCopyMode = $false

$wsOut.Range("A74").Value2 = 44055
$wsOut.Range("B74").Value2 = 498380
$wsOut.Range("C74").Value2 = 545262
$wsOut.Range("D74").Value2 = 83473
$wsOut.Range("E74").Value2 = 54666
$wsOut.Range("F74").Value2 = 26.485412737268749
$wsOut.Range("G74").Value2 = 131998
$wsOut.Range("H74").Value2 = 10528
$wsOut.Range("I74").Value2 = 12742
$wsOut.Range("J74").Value2 = 1127115

# ---------------------------------------------------------------------------
# dates_dx: new row 74
# ---------------------------------------------------------------------------
$wsDx.Range("A73:L73").Copy() | Out-Null
$wsDx.Range("A74:L74").PasteSpecial(-4122) | Out-Null
$wsDx.Application.CutCopyMode = $false

$wsDx.Range("A74").Value2 = 44055
$wsDx.Range("B74").Value2 = 0
$wsDx.Range("C74").Value2 = 1
$wsDx.Range("D74").Value2 = 0
$wsDx.Range("E74").Value2 = 0
$wsDx.Range("F74").Value2 = 1
$wsDx.Range("G74").Value2 = 0
$wsDx.Range("H74").Value2 = 0
$wsDx.Range("I74").Value2 = 0
$wsDx.Range("J74").Value2 = 0
$wsDx.Range("K74").Value2 = 0
$wsDx.Range("L74").Value2 = 4

# ---------------------------------------------------------------------------
# dates_sx: new row 74
# ---------------------------------------------------------------------------
$wsSx.Range("A73:N73").Copy() | Out-Null
$wsSx.Range("A74:N74").PasteSpecial(-4122) | Out-Null
$wsSx.Application.CutCopyMode = $false

$wsSx.Range("A74").Value2 = 44055
$wsSx.Range("B74").Value2 = 0
$wsSx.Range("C74").Value2 = 1
$wsSx.Range("D74").Value2 = 0
$wsSx.Range("E74").Value2 = 0
$wsSx.Range("F74").Value2 = 0
$wsSx.Range("G74").Value2 = 0
$wsSx.Range("H74").Value2 = 1
$wsSx.Range("I74").Value2 = 0
$wsSx.Range("J74").Value2 = 0
$wsSx.Range("K74").Value2 = 1
$wsSx.Range("L74").Value2 = 0
$wsSx.Range("M74").Value2 = 0
$wsSx.Range("N74").Value2 = 0

# ---------------------------------------------------------------------------
# dates_deaths: new row 74
# ---------------------------------------------------------------------------
$wsDeaths.Range("A73:J73").Copy() | Out-Null
$wsDeaths.Range("A74:J74").PasteSpecial(-4122) | Out-Null
$wsDeaths.Application.CutCopyMode = $false

$wsDeaths.Range("A74").Value2 = 44055
$wsDeaths.Range("B74").Value2 = 0
$wsDeaths.Range("C74").Value2 = 0
$wsDeaths.Range("D74").Value2 = 0
$wsDeaths.Range("E74").Value2 = 0
$wsDeaths.Range("F74").Value2 = 2
$wsDeaths.Range("G74").Value2 = 1
$wsDeaths.Range("H74").Value2 = 1
$wsDeaths.Range("I74").Value2 = 1
$wsDeaths.Range("J74").Value2 = 2

# ---------------------------------------------------------------------------
# control_obs: fill in the day-74 (column BV) observations
# ---------------------------------------------------------------------------

# BV1 holds the new date; format matches the plain date style used elsewhere
# for freshly-added date cells (out_vars!A74), not the shaded header style.
$wsOut.Range("A74").Copy() | Out-Null
$wsControl.Range("BV1").PasteSpecial(-4122) | Out-Null
$wsControl.Application.CutCopyMode = $false
$wsControl.Range("BV1").Value2 = 44055

$wsControl.Range("BU2").Copy() | Out-Null
$wsControl.Range("BV2").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV2").Value2 = 5250

$wsControl.Range("BU3").Copy() | Out-Null
$wsControl.Range("BV3").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV3").Value2 = 5049

$wsControl.Range("BU4").Copy() | Out-Null
$wsControl.Range("BV4").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV4").Value2 = 5049

$wsControl.Range("BU5").Copy() | Out-Null
$wsControl.Range("BV5").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV5").Value2 = 5049

$wsControl.Range("BU6").Copy() | Out-Null
$wsControl.Range("BV6").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV6").Value2 = 5049

$wsControl.Range("BU7").Copy() | Out-Null
$wsControl.Range("BV7").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV7").Value2 = 4359

$wsControl.Range("BU8").Copy() | Out-Null
$wsControl.Range("BV8").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV8").Value2 = 6990

# Row 9 separator already has a formatted (empty) BV9 cell - leave as-is.

# BV10 uses the regular style (row 10's BU cell carries a one-off highlight
# that shouldn't be propagated), so borrow the format from the row below.
$wsControl.Range("BU11").Copy() | Out-Null
$wsControl.Range("BV10").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV10").Value2 = 213

$wsControl.Range("BU11").Copy() | Out-Null
$wsControl.Range("BV11").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV11").Value2 = 213

$wsControl.Range("BU12").Copy() | Out-Null
$wsControl.Range("BV12").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV12").Value2 = 213

$wsControl.Range("BU13").Copy() | Out-Null
$wsControl.Range("BV13").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV13").Value2 = 213

$wsControl.Range("BU14").Copy() | Out-Null
$wsControl.Range("BV14").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV14").Value2 = 213

$wsControl.Range("BU15").Copy() | Out-Null
$wsControl.Range("BV15").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV15").Value2 = 148

$wsControl.Range("BU16").Copy() | Out-Null
$wsControl.Range("BV16").PasteSpecial(-4122) | Out-Null
$wsControl.Range("BV16").Value2 = 225

# Row 17 separator already has a formatted (empty) BV17 cell - leave as-is.

# BV18 is a brand-new cell (row 18 previously stopped at BU18).
$wsControl.Range("BH18").Copy() | Out-Null
$wsControl.Range("BV18").PasteSpecial(-4122) | Out-Null
$wsControl.Application.CutCopyMode = $false
$wsControl.Range("BV18").Value2 = 1220

# BV20 extends the running SUM() total one more column.
$wsControl.Range("BU20").Copy() | Out-Null
$wsControl.Range("BV20").PasteSpecial(-4122) | Out-Null
$wsControl.Application.CutCopyMode = $false
$wsControl.Range("BV20").Formula = "=SUM(BV2:BV18)"

# ---------------------------------------------------------------------------
# View state: selections move forward with the new data, and control_obs
# becomes the active sheet/tab (it was the sheet being worked on).
# ---------------------------------------------------------------------------
$wsOut.Activate()
$wsOut.Range("B75").Select() | Out-Null

$wsDx.Activate()
$wsDx.Range("L74").Select() | Out-Null

$wsSx.Activate()
$wsSx.Range("O74").Select() | Out-Null

$wsDeaths.Activate()
$wsDeaths.Range("K74").Select() | Out-Null

$wsControl.Activate()
$wsControl.Range("BW25").Select() | Out-Null
